$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the main PWM input parameters (Period, Duty-Cycle, Phase, Offset)
$ws.Range("B4").Value = 0.00005
$ws.Range("C4").Value = 50
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 120

# Fix the Duty-Cycle real (%) row formulas to use an absolute reference to $C$4
$ws.Range("C10").Formula = '=IF(((C8*($C$4/100))+C9)<=(2^$G$4-1),((C8*($C$4/100))+C9),"Out of range")'
$ws.Range("D10").Formula = '=IF(((D8*($C$4/100))+D9)<=(2^$G$4-1),((D8*($C$4/100))+D9),"Out of range")'
$ws.Range("E10").Formula = '=IF(((E8*($C$4/100))+E9)<=(2^$G$4-1),((E8*($C$4/100))+E9),"Out of range")'
$ws.Range("F10").Formula = '=IF(((F8*($C$4/100))+F9)<=(2^$G$4-1),((F8*($C$4/100))+F9),"Out of range")'
$ws.Range("G10").Formula = '=IF(((G8*($C$4/100))+G9)<=(2^$G$4-1),((G8*($C$4/100))+G9),"Out of range")'

# Move the active selection to H17 like the saved workbook
$ws.Range("H17").Select()
